# "finalização dos links da movimentação"
#
# - Rename the "calculos" sheet to "Morador".
# - Make it the active (selected) tab, with F12 as the active cell.
# - The previously-active "Contas" tab keeps its own selection (C7) but is
#   no longer the active tab.
# - Nudge the workbook window position/size to match the final saved view.

$wb = $excel.ActiveWorkbook

$wsMorador = $wb.Worksheets.Item("calculos")
$wsMorador.Name = "Morador"

# Activating this sheet and selecting a cell on it makes it the workbook's
# active tab and updates its saved selection; the sheet that used to be
# active (Contas) automatically stops being tabSelected.
$wsMorador.Activate()
$wsMorador.Range("F12").Select()

# Match the saved window geometry.
$win = $excel.ActiveWindow
$win.Left = 30
$win.Top = 840
$win.Width = 28800
$win.Height = 11295
